# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with freshly scraped figures.
#
# For Price cells whose new text parses as a plain number (e.g. "566.66",
# "1.00", "0.0645"), Excel's COM Value setter would silently convert the
# string to a double and drop the original text formatting (trailing
# zeros, etc). To keep these as literal text -- matching the source data,
# which stores every Price/Volume cell as a string -- NumberFormat is
# forced to "@" (Text) right before the assignment, then the cell style is
# restored to "Normal" afterwards so no visible/applied formatting differs
# from the original workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.927.30'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').Value = '3.128.46'
$ws.Range('E3').Value = '  -7.28%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.66'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.36'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.40%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.589'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.02%  '
$ws.Range('D9').Value = '3.127.79'
$ws.Range('E9').Value = '  -7.20%  '
$ws.Range('E10').Value = '  -6.68%  '
$ws.Range('E11').Value = '  -5.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.386'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.71%  '
$ws.Range('D13').Value = '3.661.29'
$ws.Range('E13').Value = '  -7.55%  '
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.72'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Value = '64.906.05'
$ws.Range('E16').Value = '  -1.57%  '
$ws.Range('E17').Value = '  -6.30%  '
$ws.Range('D18').Value = '3.121.02'
$ws.Range('E18').Value = '  -7.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.62'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -7.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '353.78'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.16'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.42'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.65%  '
$ws.Range('E25').Value = '  -7.89%  '
$ws.Range('D26').Value = '3.252.90'
$ws.Range('E27').Value = '  -10.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.52'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.29%  '
$ws.Range('E29').Value = '  -2.39%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('E32').Value = '  -4.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.55'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.21'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -9.06%  '
$ws.Range('E35').Value = '  -7.03%  '
$ws.Range('E36').Value = '  -5.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '159.07'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('E38').Value = '  -6.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.819'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '25.87'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.51%  '
$ws.Range('E41').Value = '  -2.50%  '
$ws.Range('D42').Value = '2.618.88'
$ws.Range('E42').Value = '  -2.48%  '
$ws.Range('E43').Value = '  -7.51%  '
$ws.Range('E44').Value = '  -4.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '39.37'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.11'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0645'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.43'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.74%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '319.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.72%  '
$ws.Range('E50').Value = '  -4.63%  '
$ws.Range('E51').Value = '  -3.76%  '
